# Append 12 new data rows (rows 206-217) to Sheet1, extending the
# normalized X-data series used for machine learning from A1:B205 to
# A1:B217 (12 more months -> A index values 204-215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A (index) / B (value) pairs to append, starting at row 206.
# B values are cast from their literal decimal-string form so the exact
# IEEE-754 double (scientific notation) is reproduced without precision
# loss.
$aValues = @(204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)
$bValues = @(
    [double]"6.106226635438361E-16",
    [double]"6.459479416000911E-15",
    [double]"3.552713678800501E-16",
    [double]"-2.467162276944792E-17",
    [double]"-8.326672684688674E-17",
    [double]"1.586032892321652E-16",
    [double]"2.960594732333751E-16",
    [double]"2.664535259100376E-16",
    [double]"2.220446049250313E-16",
    [double]"0",
    [double]"0",
    [double]"0"
)

$startRow = 206

for ($i = 0; $i -lt $aValues.Count; $i++) {
    $r = $startRow + $i

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $aValues[$i]
    # Match the formatting of the preceding A-column cells (bold,
    # bordered, centered "s=1" style) by copying the style from the row
    # directly above the newly appended one.
    $ws.Cells.Item($r - 1, 1).Copy()
    $aCell.PasteSpecial(-4122)

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $bValues[$i]
}

$excel.CutCopyMode = $false
